# Apply "more work towards final product":
#   - fill in the new carrier-word column (D) for the practice/generic pairs
#   - record the pair_kind (J) for the generic pairs that got "unique" video/audio treatment
#   - populate rows 14-21 (previously only had the row number) with their
#     kind (C) and carrier (D) values, mirroring the new "unique_video" /
#     "unique_audio" pair kinds introduced in rows 6-9
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5): new carrier word column D
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows (6-9): new pair_kind column J ("unique_video" / "unique_audio")
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: newly populated kind (C) + carrier (D) pairs
$uniqueRows = @{
    14 = @("unique_video", "look")
    15 = @("unique_video", "look")
    16 = @("unique_video", "where")
    17 = @("unique_video", "where")
    18 = @("unique_audio", "can")
    19 = @("unique_audio", "can")
    20 = @("unique_audio", "do")
    21 = @("unique_audio", "do")
}

foreach ($r in $uniqueRows.Keys) {
    $kind = $uniqueRows[$r][0]
    $carrier = $uniqueRows[$r][1]
    $ws.Cells.Item($r, 3).Value = $kind
    $ws.Cells.Item($r, 4).Value = $carrier
}
